# paper -- adding GWAS results
# Sheet1 already has the header row (observation / conclusion / extension);
# append the brainstorm rows of observations (col A) and matching
# conclusions (col B) below it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values are written in this specific order (not strict row order)
# so that the workbook's shared-string table is built up in the same
# sequence the original author typed them.
$ws.Cells.Item(3, 1).Value = "domesticated hosts show more variable lesion size"
$ws.Cells.Item(2, 1).Value = "overlap between wild and domesticated host lesion size"
$ws.Cells.Item(4, 1).Value = "isolate order changes depending on host"
$ws.Cells.Item(5, 1).Value = "bo5.10 is an average isolate"
$ws.Cells.Item(6, 1).Value = "many genes are found w large lesion effects"
$ws.Cells.Item(7, 1).Value = "effect size of largest genes depends on plant host"
$ws.Cells.Item(6, 2).Value = "strongly multigenic trait"
$ws.Cells.Item(2, 2).Value = "domestication has a host-dependent effect on resistance"
$ws.Cells.Item(3, 2).Value = "domestication did not reduce genetic diversity for lesion size"
$ws.Cells.Item(8, 1).Value = "DmWoD hits nonoverlapping with D or W"

# Widen columns A & B to fit the new (long) text, mirroring the
# "best fit" column widths recorded for the edited sheet.
$ws.Columns.Item(1).ColumnWidth = 51.3
$ws.Columns.Item(2).ColumnWidth = 55.65

# Final selection left on B8, matching the saved view state.
$ws.Range("B8").Select()
